$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.894657511337913
$ws.Range("C2").Value = 0.5413759559079826
$ws.Range("D2").Value = 0.7993492881779949
$ws.Range("E2").Value = -0.3415364410993911
$ws.Range("F2").Value = 1.380772675946646
$ws.Range("G2").Value = 1.551572196051568
$ws.Range("H2").Value = 1.523629121114938
# Row 3
$ws.Range("B3").Value = 0.4868829909126217
$ws.Range("C3").Value = 0.744856323182634
$ws.Range("D3").Value = -0.396029406094752
$ws.Range("E3").Value = 1.326279710951285
$ws.Range("F3").Value = 1.497079231056207
$ws.Range("G3").Value = 1.469136156119577
# Row 4
$ws.Range("B4").Value = -0.04030917092897979
$ws.Range("C4").Value = -1.181194900206366
$ws.Range("D4").Value = 0.5411142168396716
$ws.Range("E4").Value = 0.7119137369445936
$ws.Range("F4").Value = 0.6839706620079635
$ws.Range("G4").Value = -0.3740291588214022
$ws.Range("H4").Value = 0.7741505551085035
$ws.Range("I4").Value = 0.1041667111110185
$ws.Range("J4").Value = -0.1770340777201527
# Row 5
$ws.Range("B5").Value = -0.9071684805670841
$ws.Range("C5").Value = 0.8151406364789533
$ws.Range("D5").Value = 0.9859401565838752
$ws.Range("E5").Value = 0.9579970816472452
$ws.Range("F5").Value = -0.1000027391821206
$ws.Range("G5").Value = 1.048176974747785
$ws.Range("H5").Value = 0.3781931307503001
$ws.Range("I5").Value = 0.09699234191912887
# Row 6
$ws.Range("B6").Value = 0.873601025511372
$ws.Range("C6").Value = 1.044400545616294
$ws.Range("D6").Value = 1.016457470679664
$ws.Range("E6").Value = -0.0415423501497019
$ws.Range("F6").Value = 1.106637363780204
$ws.Range("G6").Value = 0.4366535197827188
$ws.Range("H6").Value = 0.1554527309515476
# Row 7
$ws.Range("B7").Value = 0.8572631269417244
$ws.Range("C7").Value = 0.8293200520050943
$ws.Range("D7").Value = -0.2286797688242714
$ws.Range("E7").Value = 0.9194999451056343
$ws.Range("F7").Value = 0.2495161011081493
$ws.Range("G7").Value = -0.03168468772302191
# Row 8
$ws.Range("B8").Value = 0.7156754319170088
$ws.Range("C8").Value = -0.3423243889123569
$ws.Range("D8").Value = 0.8058553250175488
$ws.Range("E8").Value = 0.1358714810200637
$ws.Range("F8").Value = -0.1453293078111075
$ws.Range("G8").Value = -0.2845013842505723
$ws.Range("H8").Value = -0.5045789468551448
$ws.Range("I8").Value = -0.3655650551191031
# Row 9
$ws.Range("B9").Value = -0.5709664978720881
$ws.Range("C9").Value = 0.5772132160578176
$ws.Range("D9").Value = -0.09277062793966739
$ws.Range("E9").Value = -0.3739714167708386
$ws.Range("F9").Value = -0.5131434932103034
$ws.Range("G9").Value = -0.733221055814876
$ws.Range("H9").Value = -0.5942071640788342
# Row 10
$ws.Range("B10").Value = 0.6742151266834819
$ws.Range("C10").Value = 0.004231282685996929
$ws.Range("D10").Value = -0.2769695061451743
$ws.Range("E10").Value = -0.4161415825846391
$ws.Range("F10").Value = -0.6362191451892116
$ws.Range("G10").Value = -0.4972052534531699
# Row 11
$ws.Range("B11").Value = 0.2471364659553615
$ws.Range("C11").Value = -0.03406432287580974
$ws.Range("D11").Value = -0.1732363993152745
$ws.Range("E11").Value = -0.3933139619198471
$ws.Range("F11").Value = -0.2543000701838053
# Row 12
$ws.Range("B12").Value = -0.7116303359412484
$ws.Range("C12").Value = -0.8508024123807132
$ws.Range("D12").Value = -1.070879974985286
$ws.Range("E12").Value = -0.931866083249244
# Row 13
$ws.Range("B13").Value = -0.5776468104600809
$ws.Range("C13").Value = -0.7977243730646535
$ws.Range("D13").Value = -0.6587104813286118
# Row 14
$ws.Range("B14").Value = -0.5891525641787865
$ws.Range("C14").Value = -0.4501386724427447
# Row 15
$ws.Range("B15").Value = -0.04484482445227361
